$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the import distribution currency/exchange-rate columns for row 4
# (blank exchange rate for import distributions)
$ws.Range("L4:N4").ClearContents()
$ws.Range("O4").ClearContents()

# Update the active selection to O4
$ws.Range("O4").Select()
